# Fruta / hortaliza, semanal
# Insert two new weekly price rows for Limón (Macroferia Regional de Talca)
# above the existing data block that starts at row 1113, shifting the
# remaining 1113..1177 rows down to 1115..1179.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 1113 (formats get
# carried over from the surrounding rows automatically, same as Excel's
# normal "insert row" behaviour).
$ws.Rows.Item(1113).Insert()
$ws.Rows.Item(1113).Insert()

# --- Row 1113 -----------------------------------------------------------
$ws.Cells.Item(1113, 1).Value = 5
$ws.Cells.Item(1113, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(1113, 3).Value = "Maule"
$ws.Cells.Item(1113, 4).Value = 44706
$ws.Cells.Item(1113, 5).Value = 7
$ws.Cells.Item(1113, 6).Value = "Fruta"
$ws.Cells.Item(1113, 7).Value = 100102
$ws.Cells.Item(1113, 8).Value = "Cítricos"
$ws.Cells.Item(1113, 9).Value = 100102003
$ws.Cells.Item(1113, 10).Value = "Limón"
$ws.Cells.Item(1113, 11).Value = "Sin especificar"
$ws.Cells.Item(1113, 12).Value = "1a amarillo"
$ws.Cells.Item(1113, 13).Value = 430
$ws.Cells.Item(1113, 14).Value = 10000
$ws.Cells.Item(1113, 15).Value = 10000
$ws.Cells.Item(1113, 16).Value = 10000
$ws.Cells.Item(1113, 17).Value = "$/malla 18 kilos"
$ws.Cells.Item(1113, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(1113, 19).Value = 556
$ws.Cells.Item(1113, 20).Value = 18

# --- Row 1114 -----------------------------------------------------------
$ws.Cells.Item(1114, 1).Value = 5
$ws.Cells.Item(1114, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(1114, 3).Value = "Maule"
$ws.Cells.Item(1114, 4).Value = 44706
$ws.Cells.Item(1114, 5).Value = 7
$ws.Cells.Item(1114, 6).Value = "Fruta"
$ws.Cells.Item(1114, 7).Value = 100102
$ws.Cells.Item(1114, 8).Value = "Cítricos"
$ws.Cells.Item(1114, 9).Value = 100102003
$ws.Cells.Item(1114, 10).Value = "Limón"
$ws.Cells.Item(1114, 11).Value = "Sin especificar"
$ws.Cells.Item(1114, 12).Value = "2a amarillo"
$ws.Cells.Item(1114, 13).Value = 230
$ws.Cells.Item(1114, 14).Value = 7000
$ws.Cells.Item(1114, 15).Value = 7000
$ws.Cells.Item(1114, 16).Value = 7000
$ws.Cells.Item(1114, 17).Value = "$/malla 14 kilos"
$ws.Cells.Item(1114, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(1114, 19).Value = 500
$ws.Cells.Item(1114, 20).Value = 14
